# Trade #112 closed at 2026-02-18 00:38:15 - unknown UNKNOWN +0.000%
#
# This script applies the following to live_trading_results.xlsx:
#   - Summary: refresh aggregate metrics (capital, P&L, trade/win counters)
#   - Strategy Status: refresh HighProbConvergence strategy row
#   - All Trades: close out trade #140 (row 141) early, and append two new
#     OPEN trades (momentum #169, MarketMaking #170)
#   - momentum: append its copy of the new OPEN trade (#169)
#   - HighProbConvergence: close out its copy of trade #140 (row 16)
#   - MarketMaking: append its copy of the new OPEN trade (#170)

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    # Force plain text so Excel's autodetection doesn't turn date-looking
    # strings (e.g. "2026-02-18") into real date serials.
    $range.NumberFormat = "@"
    $range.Value = $text
}

# ---------------------------------------------------------------------
# Summary
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1499.11
$wsSummary.Range("B4").Value = 0.22
$wsSummary.Range("B6").Value = 140
$wsSummary.Range("B7").Value = 65
$wsSummary.Range("B9").Value = 46.43

# ---------------------------------------------------------------------
# Strategy Status (row 3 = HighProbConvergence)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C3").Value = 100.42
$wsStatus.Range("D3").Value = 15
$wsStatus.Range("E3").Value = 0.42
$wsStatus.Range("F3").Value = 0.42
$wsStatus.Range("G3").Value = 73.33

# ---------------------------------------------------------------------
# All Trades
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Close out row 141 (trade #140, HighProbConvergence) early.
$wsAll.Range("G141").Value = 0.07000000000000001
$wsAll.Range("H141").Value = "CLOSED"
$wsAll.Range("I141").Value = 16.6667
$wsAll.Range("J141").Value = 0.01
$wsAll.Range("K141").Value = 100.42
$wsAll.Range("L141").Value = "early_exit"
$wsAll.Range("M141").Value = 0.17

# New row 170: momentum trade #169, still OPEN.
$wsAll.Range("A170").Value = 169
Set-TextValue $wsAll.Range("B170") "2026-02-18"
$wsAll.Range("C170").Value = "00:38:08"
$wsAll.Range("D170").Value = "momentum"
$wsAll.Range("E170").Value = "DOWN"
$wsAll.Range("F170").Value = 0.06
$wsAll.Range("H170").Value = "OPEN"
$wsAll.Range("I170").Value = 0
$wsAll.Range("J170").Value = 0
$wsAll.Range("K170").Value = 99.22374292899114
$wsAll.Range("M170").Value = 0
$wsAll.Range("N170").Value = 0
$wsAll.Range("O170").Value = 0
$wsAll.Range("P170").Value = 0.9
$wsAll.Range("Q170").Value = "Downward momentum: -1.942% over 10 samples"

# New row 171: MarketMaking trade #170, still OPEN.
$wsAll.Range("A171").Value = 170
Set-TextValue $wsAll.Range("B171") "2026-02-18"
$wsAll.Range("C171").Value = "00:38:09"
$wsAll.Range("D171").Value = "MarketMaking"
$wsAll.Range("E171").Value = "DOWN"
$wsAll.Range("F171").Value = 0.07000000000000001
$wsAll.Range("H171").Value = "OPEN"
$wsAll.Range("I171").Value = 0
$wsAll.Range("J171").Value = 0
$wsAll.Range("K171").Value = 99.19858346467944
$wsAll.Range("M171").Value = 0
$wsAll.Range("N171").Value = 0
$wsAll.Range("O171").Value = 0
$wsAll.Range("P171").Value = 0.6
$wsAll.Range("Q171").Value = "Normal spread capture: 198 bps"

# ---------------------------------------------------------------------
# momentum sheet: append its own copy of the new row (trade #169)
# ---------------------------------------------------------------------
$wsMomentum = $wb.Worksheets.Item("momentum")
$wsMomentum.Range("A44").Value = 169
Set-TextValue $wsMomentum.Range("B44") "2026-02-18"
$wsMomentum.Range("C44").Value = "00:38:08"
$wsMomentum.Range("D44").Value = "momentum"
$wsMomentum.Range("E44").Value = "DOWN"
$wsMomentum.Range("F44").Value = 0.06
$wsMomentum.Range("H44").Value = "OPEN"
$wsMomentum.Range("I44").Value = 0
$wsMomentum.Range("J44").Value = 0
$wsMomentum.Range("K44").Value = 99.22374292899114
$wsMomentum.Range("L44").Value = 0
$wsMomentum.Range("M44").Value = 0
$wsMomentum.Range("N44").Value = 0.9
$wsMomentum.Range("O44").Value = "Downward momentum: -1.942% over 10 samples"
$wsMomentum.Range("Q44").Value = 0

# ---------------------------------------------------------------------
# HighProbConvergence sheet: close out its own copy (row 16)
# ---------------------------------------------------------------------
$wsHPC = $wb.Worksheets.Item("HighProbConvergence")
$wsHPC.Range("G16").Value = 0.07000000000000001
$wsHPC.Range("H16").Value = "CLOSED"
$wsHPC.Range("I16").Value = 16.6667
$wsHPC.Range("J16").Value = 0.01
$wsHPC.Range("K16").Value = 100.42
$wsHPC.Range("P16").Value = "early_exit"
$wsHPC.Range("Q16").Value = 0.17

# ---------------------------------------------------------------------
# MarketMaking sheet: append its own copy of the new row (trade #170)
# ---------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Range("A67").Value = 170
Set-TextValue $wsMM.Range("B67") "2026-02-18"
$wsMM.Range("C67").Value = "00:38:09"
$wsMM.Range("D67").Value = "MarketMaking"
$wsMM.Range("E67").Value = "DOWN"
$wsMM.Range("F67").Value = 0.07000000000000001
$wsMM.Range("H67").Value = "OPEN"
$wsMM.Range("I67").Value = 0
$wsMM.Range("J67").Value = 0
$wsMM.Range("K67").Value = 99.19858346467944
$wsMM.Range("L67").Value = 0
$wsMM.Range("M67").Value = 0
$wsMM.Range("N67").Value = 0.6
$wsMM.Range("O67").Value = "Normal spread capture: 198 bps"
$wsMM.Range("Q67").Value = 0
